$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header E1: "fakultas" -> "prodi"
$ws.Range("E1").Value = "prodi"

# Trim leading tab/whitespace from the "prodi" data values in column E (rows 2-17)
for ($r = 2; $r -le 17; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $current = $cell.Value()
    $cell.Value = $current.Trim()
}

# Update the active selection to D8 as per the diff
$ws.Range("D8").Select()
